# "Aula 26 - Finalizando services - para departamento e funcionario"
#
# Adds the new log entry for "27. RESUMO" as row 16 of the tracking sheet:
#   B16 = aula number (27)
#   C16 = sessão       ("5. Camada de Serviço" - same session as row 15)
#   D16 = nome da aula ("27. RESUMO")
#   E16 = observação   (long note about the @Transactional resumo)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B16").Value = 27
$ws.Range("C16").Value = "5. Camada de Serviço"
$ws.Range("D16").Value = "27. RESUMO"
$ws.Range("E16").Value = "Não tem video porém tem um resumo interessante com mais detalhes sobre a anotação @Transactional e exemplos de uso, vale a pena a leitura"

# Column E uses wrap text throughout the sheet; match it for the new row so
# the longer note wraps instead of overflowing.
$ws.Range("E16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 30

# Bring the new row into view, same as the author scrolling down to it.
$ws.Range("E16").Select()
